# Add a new "booleans" worksheet between "dates" and "sharedstrings",
# containing a header row with an intentional empty cell (D3) flanked by
# non-empty headers in C3 and E3, to exercise "reading headers with empty
# cells".

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# Insert the new sheet right after "dates" (i.e. before "sharedstrings").
$datesSheet = $sheets.Item("dates")
$newSheet = $sheets.Add([System.Reflection.Missing]::Value, $datesSheet)
$newSheet.Name = "booleans"

# Header row (row 3). D3 is intentionally left empty.
$newSheet.Range("C3").Value = "boolcol"

# Boolean column C and string column D, rows 4-8.
$newSheet.Range("C4").Value = $true
$newSheet.Range("D4").Value = "a"

$newSheet.Range("C5").Value = $false
$newSheet.Range("D5").Value = "b"

$newSheet.Range("C6").Value = $true
$newSheet.Range("D6").Value = "c"

$newSheet.Range("C7").Value = $false
$newSheet.Range("D7").Value = "d"

$newSheet.Range("C8").Value = $false
$newSheet.Range("D8").Value = "e"

# Second header + numeric data column E (added last so the shared-string
# table ends up with "moredata" appended after "a".."e").
$newSheet.Range("E3").Value = "moredata"
$newSheet.Range("E4").Value = 1
$newSheet.Range("E5").Value = 2
$newSheet.Range("E6").Value = 3
$newSheet.Range("E7").Value = 4
$newSheet.Range("E8").Value = 5

# Make this the active sheet/selection, matching the saved view state.
$newSheet.Range("D4").Select() | Out-Null
